$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data point was added at the top of the data table (row 8).
# All the existing data rows from 8 downward shift down by one row.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new record's data.
$ws.Cells.Item(8, 1).Value = 9
$ws.Cells.Item(8, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(8, 3).Value = "Metropolitana"
$ws.Cells.Item(8, 4).Value = 44503
$ws.Cells.Item(8, 5).Value = 13
$ws.Cells.Item(8, 6).Value = 100112005
$ws.Cells.Item(8, 7).Value = "Puerro"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 97
$ws.Cells.Item(8, 11).Value = 7000
$ws.Cells.Item(8, 12).Value = 8000
$ws.Cells.Item(8, 13).Value = 7505
$ws.Cells.Item(8, 14).Value = "`$/paquete 20 unidades"
$ws.Cells.Item(8, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(8, 16).Value = 375
$ws.Cells.Item(8, 17).Value = 20
$ws.Cells.Item(8, 18).Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of the table.
$ws.Cells.Item(8, 4).NumberFormat = $ws.Cells.Item(9, 4).NumberFormat
